$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows for Aug 9, 10, 11 2022 (serials 44782-44784)
$data = @(
    @(44782, 0, 0, 0),
    @(44783, 0, 1, 0),
    @(44784, 0, 0, 0)
)

$r = 27
foreach ($row in $data) {
    # Row height must be set before the cell writes to stick (matches the
    # 13.8pt height already used by the preceding rows 24-26)
    $ws.Rows.Item($r).RowHeight = 13.8

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    # Match formatting of the preceding rows (date column uses M/D/YYYY,
    # the rest use General number format)
    $ws.Cells.Item($r, 1).NumberFormat = "M/D/YYYY"
    $ws.Cells.Item($r, 2).NumberFormat = "General"
    $ws.Cells.Item($r, 3).NumberFormat = "General"
    $ws.Cells.Item($r, 4).NumberFormat = "General"

    $r++
}

$ws.Range("C29").Select()
